$wb = $excel.ActiveWorkbook

$wsWater = $wb.Worksheets.Item("EZ Water Adjustment")
$wsRaw   = $wb.Worksheets.Item("Raw Text Format")

# --- EZ Water Adjustment sheet -------------------------------------------
$wsWater.Activate()

# Bump the total water volume (liters) input - this cascades through every
# downstream formula (strike/sparge ratios, mineral additions, sulfate /
# chloride ratio, etc.) via normal recalculation.
$wsWater.Range("D8").Value = 14.4

# Recorded view state: scrolled down a bit, with D8 selected.
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$wsWater.Range("D8").Select()

# --- Raw Text Format sheet -------------------------------------------------
$wsRaw.Activate()
$wsRaw.Range("B18").Select()

# Leave the EZ Water Adjustment sheet as the active/selected tab, matching
# the saved workbook state.
$wsWater.Activate()
